# Adds the new "tween" (easing function) reference sheet, and adds a
# "level" column to the "weapon" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New "tween" sheet, appended after the last existing sheet ("weapon").
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$tween = $wb.Worksheets.Add($null, $lastSheet)
$tween.Name = "tween"

$tweenValues = @(
    "id",
    "string",
    "缓动函数",
    "linear",
    "smooth",
    "fade",
    "constant",
    "quadIn",
    "quadOut",
    "quadInOut",
    "quadOutIn",
    "cubicIn",
    "cubicOut",
    "cubicInOut",
    "cubicOutIn",
    "quartIn",
    "quartOut",
    "quartInOut",
    "quartOutIn",
    "quintIn",
    "quintOut",
    "quintInOut",
    "quintOutIn",
    "sineIn",
    "sineOut",
    "sineInOut",
    "sineOutIn",
    "expoIn",
    "expoOut",
    "expoInOut",
    "expoOutIn",
    "circIn",
    "circOut",
    "circInOut",
    "circOutIn",
    "elasticIn",
    "elasticOut",
    "elasticInOut",
    "elasticOutIn",
    "backIn",
    "backOut",
    "backInOut",
    "backOutIn",
    "bounceIn",
    "bounceOut",
    "bounceInOut",
    "bounceOutIn",
)

for ($i = 0; $i -lt $tweenValues.Length; $i++) {
    $tween.Cells.Item($i + 1, 1).Value = $tweenValues[$i]
}

$tween.Columns("A").ColumnWidth = 23.33203125

# ---------------------------------------------------------------------
# 2) "weapon" sheet: insert a new "level" column right after "weaponId".
# ---------------------------------------------------------------------
$weapon = $wb.Worksheets.Item("weapon")
$weapon.Columns("B").Insert()

$weapon.Range("B1").Value = "level"
$weapon.Range("B2").Value = "int"
$weapon.Range("B3").Value = "等级"

$weaponLevels = @(1, 2, 3, 4, 5, 6, 7)
for ($i = 0; $i -lt $weaponLevels.Length; $i++) {
    $weapon.Cells.Item($i + 4, 2).Value = $weaponLevels[$i]
}

# ---------------------------------------------------------------------
# 3) Restore per-sheet selections.
# ---------------------------------------------------------------------
$activity = $wb.Worksheets.Item("activity")
$activity.Range("D14").Select()

$weapon.Range("E15").Select()

# Activate "tween" last so it becomes the active/selected tab, matching
# the workbook's new activeTab.
$tween.Range("E20").Select()
